# Add two more sheets ("Sheet2", "Sheet3") that are copies of Sheet1's
# data/formatting, positioned after Sheet1, then adjust the view/selection
# state on all three sheets to match the target layout.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheets right after Sheet1 / Sheet2 so the final tab order
# is Sheet1, Sheet2, Sheet3 (Worksheets.Add() with no placement args would
# otherwise insert before the active sheet).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws3 = $wb.Worksheets.Add($null, $ws2)

# Populate the new sheets with the same data + styles as Sheet1.
$ws1.Range("A1:B37").Copy($ws2.Range("A1:B37"))
$ws1.Range("A1:B37").Copy($ws3.Range("A1:B37"))

# Sheet1 + Sheet2: selection becomes the whole used range.
$ws1.Range("A1:B37").Select()
$ws2.Range("A1:B37").Select()

# Sheet3: selection is a single cell, D10 - and it's the tab left active.
$ws3.Range("D10").Select()
$ws3.Activate()
